$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''41.853.29'
$ws.Cells.Item(2, 5).Value = '  +4.22%  '
$ws.Cells.Item(3, 4).Value = '''2.276.02'
$ws.Cells.Item(3, 5).Value = '  +2.26%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).Value = '''305.83'
$ws.Cells.Item(5, 5).Value = '  +4.15%  '
$ws.Cells.Item(6, 4).Value = '''92.64'
$ws.Cells.Item(6, 5).Value = '  +5.37%  '
$ws.Cells.Item(7, 5).Value = '  +3.84%  '
$ws.Cells.Item(9, 4).Value = '''0.485'
$ws.Cells.Item(9, 5).Value = '  +3.32%  '
$ws.Cells.Item(10, 4).Value = '''32.80'
$ws.Cells.Item(10, 5).Value = '  +6.85%  '
$ws.Cells.Item(11, 4).Value = '''53.71'
$ws.Cells.Item(11, 5).Value = '  +5.59%  '
$ws.Cells.Item(12, 5).Value = '  +2.17%  '
$ws.Cells.Item(13, 5).Value = '  +1.46%  '
$ws.Cells.Item(14, 5).Value = '  +3.63%  '
$ws.Cells.Item(15, 4).Value = '''2.627.88'
$ws.Cells.Item(15, 5).Value = '  +2.19%  '
$ws.Cells.Item(16, 5).Value = '  +3.44%  '
$ws.Cells.Item(17, 4).Value = '''2.266.54'
$ws.Cells.Item(17, 5).Value = '  +0.64%  '
$ws.Cells.Item(18, 5).Value = '  +3.75%  '
$ws.Cells.Item(19, 4).Value = '''41.795.83'
$ws.Cells.Item(19, 5).Value = '  +4.18%  '
$ws.Cells.Item(20, 4).Value = '''12.39'
$ws.Cells.Item(20, 5).Value = '  +10.12%  '
$ws.Cells.Item(21, 5).Value = '  +2.23%  '
$ws.Cells.Item(22, 4).Value = '''5.94'
$ws.Cells.Item(22, 5).Value = '  +2.74%  '
$ws.Cells.Item(23, 4).Value = '''67.13'
$ws.Cells.Item(23, 5).Value = '  +2.22%  '
$ws.Cells.Item(24, 4).Value = '''242.71'
$ws.Cells.Item(24, 5).Value = '  +2.84%  '
$ws.Cells.Item(25, 5).Value = '  +5.02%  '
$ws.Cells.Item(26, 5).Value = '  +0.15%  '
$ws.Cells.Item(27, 5).Value = '  +5.37%  '
$ws.Cells.Item(28, 4).Value = '''24.32'
$ws.Cells.Item(28, 5).Value = '  +4.75%  '
$ws.Cells.Item(29, 4).Value = '''9.62'
$ws.Cells.Item(29, 5).Value = '  +3.07%  '
$ws.Cells.Item(30, 5).Value = '  +0.05%  '
$ws.Cells.Item(31, 4).Value = '''34.26'
$ws.Cells.Item(31, 5).Value = '  +7.33%  '
$ws.Cells.Item(32, 4).Value = '''158.91'
$ws.Cells.Item(32, 5).Value = '  +0.11%  '
$ws.Cells.Item(33, 5).Value = '  -0.03%  '
$ws.Cells.Item(34, 5).Value = '  +4.41%  '
$ws.Cells.Item(35, 4).Value = '''0.0750'
$ws.Cells.Item(35, 5).Value = '  +4.83%  '
$ws.Cells.Item(36, 5).Value = '  +0.90%  '
$ws.Cells.Item(37, 4).Value = '''17.07'
$ws.Cells.Item(37, 5).Value = '  +8.89%  '
$ws.Cells.Item(38, 5).Value = '  +2.02%  '
$ws.Cells.Item(39, 5).Value = '  +2.75%  '
$ws.Cells.Item(40, 4).Value = '''0.104'
$ws.Cells.Item(40, 5).Value = '  +5.08%  '
$ws.Cells.Item(41, 5).Value = '  +3.17%  '
$ws.Cells.Item(42, 5).Value = '  +4.52%  '
$ws.Cells.Item(43, 4).Value = '''2.072.43'
$ws.Cells.Item(43, 5).Value = '  -0.75%  '
$ws.Cells.Item(44, 5).Value = '  +1.63%  '
$ws.Cells.Item(45, 5).Value = '  +3.35%  '
$ws.Cells.Item(46, 5).Value = '  +3.11%  '
$ws.Cells.Item(47, 5).Value = '  +5.86%  '
$ws.Cells.Item(48, 5).Value = '  +7.57%  '
$ws.Cells.Item(49, 2).Value = 'BitcoinSV'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Cells.Item(49, 4).Value = '''73.11'
$ws.Cells.Item(49, 5).Value = '  +7.45%  '
$ws.Cells.Item(50, 2).Value = 'Stacks'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(50, 4).Value = '''1.52'
$ws.Cells.Item(50, 5).Value = '  +3.46%  '
$ws.Cells.Item(51, 5).Value = '  +3.63%  '
